$d = $word.ActiveDocument

# 1. Update the letter date
$d.Content.Find.Execute("September 19, 2025", $false, $false, $false, $false, $false, $true, 1, $false, "September 21, 2025", 2) | Out-Null

# 2. Split the mailing address paragraph "919 Story Road, San Jose CA 95122"
#    into two paragraphs: "919 Story Road" and "San Jose, CA 95122".
#    This is the first occurrence (the sender/recipient block), not the
#    "PROPERTY ADDRESS" table cell occurrence further down.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "919 Story Road, San Jose CA 95122*") {
        $p.Range.Text = "919 Story Road`rSan Jose, CA 95122"
        break
    }
}

# 3. Remove the empty "NoSpacing" paragraph that immediately follows the
#    "Board of Directors" signature line.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Vietnam Town Condominium Owners Association Board of Directors*") {
        $next = $p.Next()
        $next.Range.Delete()
        break
    }
}
